# Update the Fitness (column C) values for rows 2-93 (Generation 0-91)
# on the active sheet, matching the new run's logged values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    12384,12384,12163,12163,10800,10503,10503,10503,10503,10503,
    9897,9897,9275,9275,9275,9041,9041,9041,9025,8706,
    8706,8706,8706,8706,8706,8706,8345,8345,8345,8345,
    8345,8345,8345,8345,8345,8345,8345,8345,8020,8020,
    8020,8020,8020,8020,8020,8020,8020,7970,7970,7970,
    7970,7970,7970,7970,7892,7892,7892,7892,7892,7892,
    7892,7892,7892,7892,7892,7892,7892,7892,7892,7892,
    7892,7892,7892,7872,7872,7872,7872,7872,7872,7872,
    7872,7872,7872,7748,7748,7748,7748,7748,7748,7748,
    7748,7748
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
